$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 40, column A ("phone") was stored as text "71277620"; the canonical
# edit turns it into a real number, matching every other phone cell in
# the sheet.
$ws.Range("A40").Value = 71277620

# Append the new payment row 41: 71277620 (Cash) 2025-08-18T17:29:26
# The phone number must stay TEXT (it keeps its leading format as typed
# into the source system), so force text entry via a leading apostrophe
# and then drop back to the workbook's default "Normal" style so it
# doesn't pick up a quote-prefix style like a manually-typed cell would.
$ws.Range("A41").Value = "'71277620"
$ws.Range("A41").Style = "Normal"

$ws.Range("B41").Value = ""
$ws.Range("C41").Value = "Cash"
$ws.Range("D41").Value = "2025-08-18T17:29:26"
$ws.Range("E41").Value = 100
$ws.Range("F41").Value = ""
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 100
